# Scheduled-runner refresh: updates computed market-price / profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the Excalibur_Profits sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 942.95
$ws.Range("I28").Value = 1027
$ws.Range("J28").Value = 816.875
$ws.Range("K28").Value = 1027
$ws.Range("L28").Value = 816.875
$ws.Range("M28").Value = -542
$ws.Range("N28").Value = -1786.875

$ws.Range("H87").Value = 67534.14
$ws.Range("J87").Value = 92184.75
$ws.Range("L87").Value = 92184.75
$ws.Range("N87").Value = -94680.75

$ws.Range("H90").Value = 67534.14
$ws.Range("J90").Value = 92184.75
$ws.Range("L90").Value = 276554.25
$ws.Range("N90").Value = -289034.25

$ws.Range("H92").Value = 2212.84
$ws.Range("I92").Value = 2299.0557
$ws.Range("K92").Value = 2299.0557
$ws.Range("M92").Value = -1051.0557

$ws.Range("H107").Value = 450.29413
$ws.Range("I107").Value = 447.1875
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 447.1875
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1472.8125
$ws.Range("N107").Value = -4340

$ws.Range("H112").Value = 10217.444
$ws.Range("J112").Value = 13787.947
$ws.Range("L112").Value = 41363.841
$ws.Range("N112").Value = -43579.841

$ws.Range("H133").Value = 69194.5
$ws.Range("J133").Value = 69194.5
$ws.Range("L133").Value = 69194.5
$ws.Range("N133").Value = -79314.5

$ws.Range("H136").Value = 96499.5
$ws.Range("J136").Value = 96499.5
$ws.Range("L136").Value = 96499.5
$ws.Range("N136").Value = -106699.5

$ws.Range("H137").Value = 2008
$ws.Range("I137").Value = 1686
$ws.Range("K137").Value = 5058
$ws.Range("M137").Value = -2508

$ws.Range("H138").Value = 1692.525
$ws.Range("I138").Value = 1137.0646
$ws.Range("J138").Value = 3605.7778
$ws.Range("K138").Value = 3411.1938
$ws.Range("L138").Value = 10817.3334
$ws.Range("M138").Value = 1728.8062
$ws.Range("N138").Value = -21097.3334

$ws.Range("H141").Value = 1885.7142
$ws.Range("I141").Value = 1839
$ws.Range("K141").Value = 5517
$ws.Range("M141").Value = -337

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3211.5334
$ws.Range("I45").Value = 3752.2727
$ws.Range("K45").Value = 3752.2727
$ws.Range("M45").Value = -3375.2727

$ws.Range("H61").Value = 5339.7295
$ws.Range("I61").Value = 5325.148
$ws.Range("K61").Value = 5325.148
$ws.Range("M61").Value = -5113.148

$ws.Range("H110").Value = 1000.25
$ws.Range("I110").Value = 1000.25
$ws.Range("K110").Value = 1000.25
$ws.Range("M110").Value = 1044.75

$ws.Range("H134").Value = 71389.8
$ws.Range("J134").Value = 71389.8
$ws.Range("L134").Value = 71389.8
$ws.Range("N134").Value = -81529.8

$ws.Range("H136").Value = 5339.7295
$ws.Range("I136").Value = 5325.148
$ws.Range("K136").Value = 15975.444
$ws.Range("M136").Value = -13425.444

$ws.Range("H141").Value = 65992.664
$ws.Range("J141").Value = 65992.664
$ws.Range("L141").Value = 65992.664
$ws.Range("N141").Value = -76352.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2842.5806
$ws.Range("I134").Value = 1405.4584
$ws.Range("K134").Value = 4216.3752
$ws.Range("M134").Value = -1681.3752

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3899.6
$ws.Range("I22").Value = 433
$ws.Range("K22").Value = 433
$ws.Range("M22").Value = -83

$ws.Range("H58").Value = 3076.389
$ws.Range("I58").Value = 1065.5
$ws.Range("K58").Value = 1065.5
$ws.Range("M58").Value = -862.5

$ws.Range("H68").Value = 85594.5
$ws.Range("J68").Value = 85594.5
$ws.Range("L68").Value = 85594.5
$ws.Range("N68").Value = -87092.5

$ws.Range("H71").Value = 85594.5
$ws.Range("J71").Value = 85594.5
$ws.Range("L71").Value = 256783.5
$ws.Range("N71").Value = -264271.5

$ws.Range("H132").Value = 4810954.5
$ws.Range("I132").Value = 3379.38
$ws.Range("J132").Value = 125000340
$ws.Range("K132").Value = 10138.14
$ws.Range("L132").Value = 375001020
$ws.Range("M132").Value = -7608.139999999999
$ws.Range("N132").Value = -375006080

$ws.Range("H136").Value = 3076.389
$ws.Range("I136").Value = 1065.5
$ws.Range("K136").Value = 3196.5
$ws.Range("M136").Value = -646.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 180.64706
$ws.Range("J17").Value = 587.25
$ws.Range("L17").Value = 1761.75
$ws.Range("N17").Value = -2099.75

$ws.Range("H56").Value = 7924.8423
$ws.Range("I56").Value = 7924.8423
$ws.Range("K56").Value = 7924.8423
$ws.Range("M56").Value = -7394.8423

$ws.Range("H70").Value = 3203.2222
$ws.Range("I70").Value = 1340.6666
$ws.Range("J70").Value = 6928.3335
$ws.Range("K70").Value = 4021.9998
$ws.Range("L70").Value = 20785.0005
$ws.Range("M70").Value = -3706.9998
$ws.Range("N70").Value = -21415.0005

$ws.Range("H73").Value = 3203.2222
$ws.Range("I73").Value = 1340.6666
$ws.Range("J73").Value = 6928.3335
$ws.Range("K73").Value = 4021.9998
$ws.Range("L73").Value = 20785.0005
$ws.Range("M73").Value = -2929.9998
$ws.Range("N73").Value = -22969.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2907.4736
$ws.Range("I113").Value = 2401.1538
$ws.Range("J113").Value = 4004.5
$ws.Range("K113").Value = 2401.1538
$ws.Range("L113").Value = 4004.5
$ws.Range("M113").Value = -231.1538
$ws.Range("N113").Value = -8344.5

$ws.Range("H132").Value = 593971.6
$ws.Range("I132").Value = 938032
$ws.Range("K132").Value = 2814096
$ws.Range("M132").Value = -2811566

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1050.25

$ws.Range("H132").Value = 770623.0600000001
$ws.Range("I132").Value = 1082676
$ws.Range("K132").Value = 3248028
$ws.Range("M132").Value = -3245498

$ws.Range("H136").Value = 2991.543
$ws.Range("I136").Value = 2342.5757
$ws.Range("J136").Value = 13699.5
$ws.Range("K136").Value = 7027.7271
$ws.Range("L136").Value = 41098.5
$ws.Range("M136").Value = -4477.7271
$ws.Range("N136").Value = -46198.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 681966
$ws.Range("I132").Value = 846249.25
$ws.Range("K132").Value = 2538747.75
$ws.Range("M132").Value = -2536217.75

$ws.Range("H136").Value = 12691785
$ws.Range("I136").Value = 14617075
$ws.Range("K136").Value = 43851225
$ws.Range("M136").Value = -43848675
